# Regenerate merged AHB files
#
# The "Änderung" marker (shared string 238, style index 7 = bold / gold
# font, centered, gray fill) is removed from column L for rows 76-108
# (except the already-blank row 79), leaving an empty, centered, gray-filled
# cell (style index 4) instead.
#
# Additionally, the seven "new data element" section rows (83, 86, 90, 95,
# 99, 103, 106) get their whole row's formatting switched from the plain
# "data row" look (style 5 = white fill) to the "header-ish" look used
# elsewhere in the sheet (style 2 = gray fill for most columns, style 3 =
# bold+gray fill for column B), matching the template already present in
# row 2 of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Rows where the whole row's formatting is swapped to the "section header"
# look (gray fill for everything, bold for column B), based on the existing
# template row 2.
$fullRestyleRows = @(83, 86, 90, 95, 99, 103, 106)

# Rows where only column L loses its "ÄNDERUNG" marker/style, based on the
# existing template cell L3.
$simpleRows = @(76, 77, 78, 80, 81, 82, 84, 85, 87, 88, 89, 91, 92, 93, 94, 96, 97, 98, 100, 101, 102, 104, 105, 107, 108)

# --- Full-row restyle ---------------------------------------------------
$rowTemplate = $ws.Range("A2:V2")
$rowTemplate.Copy()
foreach ($r in $fullRestyleRows) {
    $dstRow = $ws.Range("A" + $r + ":V" + $r)
    $dstRow.PasteSpecial($xlPasteFormats)
}

# --- Simple column-L restyle --------------------------------------------
$cellTemplate = $ws.Range("L3")
$cellTemplate.Copy()
foreach ($r in $simpleRows) {
    $dstCell = $ws.Range("L" + $r)
    $dstCell.PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

# --- Clear the removed "ÄNDERUNG" text from column L on every affected row
$allRows = $fullRestyleRows + $simpleRows
foreach ($r in $allRows) {
    $ws.Range("L" + $r).Value = ""
}
